# Apply "Fixed Stimulus Absolute Timestamps" edit:
#  - Rename the 5 sheets with new timestamp-based names
#  - Update the B2:B(n) cell values (filenames / labels) on each sheet

$wb = $excel.ActiveWorkbook

# --- Rename sheets (order matches workbook.xml: GNG, NB, RS, TOL, vSAT) ---
$wsGNG  = $wb.Worksheets.Item(1)
$wsNB   = $wb.Worksheets.Item(2)
$wsRS   = $wb.Worksheets.Item(3)
$wsTOL  = $wb.Worksheets.Item(4)
$wsvSAT = $wb.Worksheets.Item(5)

$wsGNG.Name  = "GNG_TO-16504778502716854"
$wsNB.Name   = "NB_TO-1650477852526959"
$wsRS.Name   = "RS_TO-16504778525279262"
$wsTOL.Name  = "TOL_TO-16504778525729594"
$wsvSAT.Name = "vSAT_TO-16504778526368186"

# --- Sheet 1: GNG_TO ---
$wsGNG.Range("B2").Value = "go_stims-1650477850236686.csv"
$wsGNG.Range("B3").Value = "GNG_stims-16504778502547188.csv"
$wsGNG.Range("B4").Value = "go_stims-16504778502566907.csv"
$wsGNG.Range("B5").Value = "GNG_stims-165047785027072.csv"

# --- Sheet 2: NB_TO ---
$wsNB.Range("B2").Value = "OB-1650477851909924.csv"
$wsNB.Range("B3").Value = "ZB-match_4-16504778506137214.csv"
$wsNB.Range("B4").Value = "OB-16504778520089262.csv"
$wsNB.Range("B5").Value = "OB-16504778518259592.csv"
$wsNB.Range("B6").Value = "TB-16504778525059602.csv"
$wsNB.Range("B7").Value = "TB-1650477852077924.csv"
$wsNB.Range("B8").Value = "ZB-match_8-16504778503026826.csv"
$wsNB.Range("B9").Value = "TB-16504778520439236.csv"
$wsNB.Range("B10").Value = "ZB-match_9-16504778503836865.csv"

# --- Sheet 3: RS_TO ---
$wsRS.Range("B2").Value = "eyes open"
$wsRS.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL_TO ---
$wsTOL.Range("B2").Value = "MM_stims-16504778525419292.csv"
$wsTOL.Range("B3").Value = "ZM_stims-16504778525289254.csv"
$wsTOL.Range("B4").Value = "MM_stims-16504778525579612.csv"
$wsTOL.Range("B5").Value = "ZM_stims-16504778525419292.csv"
$wsTOL.Range("B6").Value = "MM_stims-16504778525729594.csv"
$wsTOL.Range("B7").Value = "ZM_stims-16504778525579612.csv"

# --- Sheet 5: vSAT_TO ---
$wsvSAT.Range("B2").Value = "vSAT_stims-16504778526048186.csv"
$wsvSAT.Range("B3").Value = "SAT_stims-16504778525888197.csv"
$wsvSAT.Range("B4").Value = "SAT_stims-16504778525769293.csv"
$wsvSAT.Range("B5").Value = "vSAT_stims-16504778526208189.csv"
